$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2-25)
$bArr = New-Object "object[,]" 24,1
$bArr[0,0] = 0.7015001070189442
$bArr[1,0] = 0.6828486133654508
$bArr[2,0] = 0.6716808960933633
$bArr[3,0] = 0.6672019392284767
$bArr[4,0] = 0.6664625730276157
$bArr[5,0] = 0.6716201992517199
$bArr[6,0] = 0.6950103245349624
$bArr[7,0] = 0.7431169459979401
$bArr[8,0] = 0.7798060967122922
$bArr[9,0] = 0.7967851712134859
$bArr[10,0] = 0.8032558449861256
$bArr[11,0] = 0.8018604509625504
$bArr[12,0] = 0.7973166972077763
$bArr[13,0] = 0.7945388503914614
$bArr[14,0] = 0.7787022476541665
$bArr[15,0] = 0.7690606743891237
$bArr[16,0] = 0.7635423273934805
$bArr[17,0] = 0.7616786018139408
$bArr[18,0] = 0.7700842215780881
$bArr[19,0] = 0.7986501971940072
$bArr[20,0] = 0.8175588548762676
$bArr[21,0] = 0.8074452288713303
$bArr[22,0] = 0.7696213988517115
$bArr[23,0] = 0.7298651882808826
$ws.Range("B2:B25").Value = $bArr

# Columns D:N (rows 2-25)
$dnArr = New-Object "object[,]" 24,11
$dnArr[0,0] = 0.09636837988456293
$dnArr[0,1] = 0.1287761185195411
$dnArr[0,2] = 2.134116765214429
$dnArr[0,3] = 1.480204298445784
$dnArr[0,4] = 1.35901841967636
$dnArr[0,5] = 1.396732643407024
$dnArr[0,6] = 0.1897913559789171
$dnArr[0,7] = 0.6071898264123092
$dnArr[0,8] = 0.3082827602104743
$dnArr[0,9] = 0.2321165219618564
$dnArr[0,10] = 2.592149438367114
$dnArr[1,0] = 0.09629881893960146
$dnArr[1,1] = 0.1288882573821837
$dnArr[1,2] = 2.138013436790757
$dnArr[1,3] = 1.480297368150204
$dnArr[1,4] = 1.363992512497063
$dnArr[1,5] = 1.405896513800375
$dnArr[1,6] = 0.1895722996663167
$dnArr[1,7] = 0.5552245767987358
$dnArr[1,8] = 0.3012654240460506
$dnArr[1,9] = 0.2263363731392154
$dnArr[1,10] = 2.614493330738117
$dnArr[2,0] = 0.09626877645010445
$dnArr[2,1] = 0.1289647210967535
$dnArr[2,2] = 2.141342489945714
$dnArr[2,3] = 1.481117734135125
$dnArr[2,4] = 1.367577633699867
$dnArr[2,5] = 1.412012641346664
$dnArr[2,6] = 0.1894411974974339
$dnArr[2,7] = 0.5235165051146566
$dnArr[2,8] = 0.2970965463758262
$dnArr[2,9] = 0.2228862854765943
$dnArr[2,10] = 2.628916241808426
$dnArr[3,0] = 0.0962597334852866
$dnArr[3,1] = 0.1289977971083892
$dnArr[3,2] = 2.142934793397309
$dnArr[3,3] = 1.481643980562211
$dnArr[3,4] = 1.369172252910133
$dnArr[3,5] = 1.414628201432365
$dnArr[3,6] = 0.1893886301282528
$dnArr[3,7] = 0.5106456590079631
$dnArr[3,8] = 0.2954329624800636
$dnArr[3,9] = 0.2215053541390866
$dnArr[3,10] = 2.634970695155394
$dnArr[4,0] = 0.0962584255727652
$dnArr[4,1] = 0.1290034051872402
$dnArr[4,2] = 2.143213435062421
$dnArr[4,3] = 1.481742956826466
$dnArr[4,4] = 1.369445114292944
$dnArr[4,5] = 1.415069956605279
$dnArr[4,6] = 0.1893799532584008
$dnArr[4,7] = 0.5085115254981929
$dnArr[4,8] = 0.2951588591116376
$dnArr[4,9] = 0.2212775656499915
$dnArr[4,10] = 2.635986724580768
$dnArr[5,0] = 0.09626864151906922
$dnArr[5,1] = 0.1289651594075989
$dnArr[5,2] = 2.141363009753078
$dnArr[5,3] = 1.481124054090117
$dnArr[5,4] = 1.367598597969163
$dnArr[5,5] = 1.412047416779444
$dnArr[5,6] = 0.1894404850775331
$dnArr[5,7] = 0.5233427195002776
$dnArr[5,8] = 0.2970739677519418
$dnArr[5,9] = 0.2228675603319701
$dnArr[5,10] = 2.628997177545198
$dnArr[6,0] = 0.09634177427574642
$dnArr[6,1] = 0.1288132062136804
$dnArr[6,2] = 2.135266087703712
$dnArr[6,3] = 1.480077991740785
$dnArr[6,4] = 1.360623350658798
$dnArr[6,5] = 1.399790819227558
$dnArr[6,6] = 0.1897151215683728
$dnArr[6,7] = 0.5892312063767804
$dnArr[6,8] = 0.3058342342705345
$dnArr[6,9] = 0.2301030636418737
$dnArr[6,10] = 2.599707601402448
$dnArr[7,0] = 0.09658513133600266
$dnArr[7,1] = 0.1285754918161366
$dnArr[7,2] = 2.130733573261367
$dnArr[7,3] = 1.484083400929734
$dnArr[7,4] = 1.35115335562142
$dnArr[7,5] = 1.379634834917567
$dnArr[7,6] = 0.1902805583271352
$dnArr[7,7] = 0.7200023486641385
$dnArr[7,8] = 0.3241184480808528
$dnArr[7,9] = 0.2450723488526876
$dnArr[7,10] = 2.547850770781498
$dnArr[8,0] = 0.0968241514780992
$dnArr[8,1] = 0.1284374373496862
$dnArr[8,2] = 2.131922380607591
$dnArr[8,3] = 1.490722594755923
$dnArr[8,4] = 1.346755837485716
$dnArr[8,5] = 1.367185018691089
$dnArr[8,6] = 0.1907122946742419
$dnArr[8,7] = 0.8170254932047669
$dnArr[8,8] = 0.3382223063951528
$dnArr[8,9] = 0.2565412637944533
$dnArr[8,10] = 2.513147371847896
$dnArr[9,0] = 0.09694582372810245
$dnArr[9,1] = 0.1283825503723044
$dnArr[9,2] = 2.133443209755583
$dnArr[9,3] = 1.494546657422177
$dnArr[9,4] = 1.345310125216969
$dnArr[9,5] = 1.362032213900214
$dnArr[9,6] = 0.190912232958496
$dnArr[9,7] = 0.8613680992391153
$dnArr[9,8] = 0.3447834888583259
$dnArr[9,9] = 0.2618601031750885
$dnArr[9,10] = 2.498096018825002
$dnArr[10,0] = 0.09699374617140322
$dnArr[10,1] = 0.1283629019570192
$dnArr[10,2] = 2.13415989126878
$dnArr[10,3] = 1.49611036347244
$dnArr[10,4] = 1.344842342036458
$dnArr[10,5] = 1.360154320590375
$dnArr[10,6] = 0.1909884513024096
$dnArr[10,7] = 0.8781888305387042
$dnArr[10,8] = 0.3472888376835641
$dnArr[10,9] = 0.2638887068645133
$dnArr[10,10] = 2.492502162729036
$dnArr[11,0] = 0.09698334324014013
$dnArr[11,1] = 0.1283670830983936
$dnArr[11,2] = 2.133999282406194
$dnArr[11,3] = 1.495768448958799
$dnArr[11,4] = 1.344939545362948
$dnArr[11,5] = 1.360555496871214
$dnArr[11,6] = 0.1909720138569924
$dnArr[11,7] = 0.8745648969103286
$dnArr[11,8] = 0.3467483440463042
$dnArr[11,9] = 0.2634511690897838
$dnArr[11,10] = 2.493702195534564
$dnArr[12,0] = 0.09694972937815294
$dnArr[12,1] = 0.1283809111276728
$dnArr[12,2] = 2.133499351067087
$dnArr[12,3] = 1.494672987202975
$dnArr[12,4] = 1.345270044110492
$dnArr[12,5] = 1.36187624861811
$dnArr[12,6] = 0.1909184933619059
$dnArr[12,7] = 0.8627513696981453
$dnArr[12,8] = 0.3449891897978148
$dnArr[12,9] = 0.2620267082798975
$dnArr[12,10] = 2.497633689136823
$dnArr[13,0] = 0.09692938017105845
$dnArr[13,1] = 0.1283895290887884
$dnArr[13,2] = 2.133211457471774
$dnArr[13,3] = 1.494017042601584
$dnArr[13,4] = 1.345482857609085
$dnArr[13,5] = 1.362694798949633
$dnArr[13,6] = 0.1908857763321166
$dnArr[13,7] = 0.8555190238859041
$dnArr[13,8] = 0.3439143585871847
$dnArr[13,9] = 0.2611560660339407
$dnArr[13,10] = 2.500055616678985
$dnArr[14,0] = 0.09681645919096127
$dnArr[14,1] = 0.1284411834310364
$dnArr[14,2] = 2.13184269869592
$dnArr[14,3] = 1.490488860893151
$dnArr[14,4] = 1.346861472893778
$dnArr[14,5] = 1.367532036511143
$dnArr[14,6] = 0.1906992992083065
$dnArr[14,7] = 0.8141317066549902
$dnArr[14,8] = 0.3377964307576633
$dnArr[14,9] = 0.2561956977989865
$dnArr[14,10] = 2.514145812478947
$dnArr[15,0] = 0.0967504909089989
$dnArr[15,1] = 0.1284748973670289
$dnArr[15,2] = 2.131253911980423
$dnArr[15,3] = 1.488530339358292
$dnArr[15,4] = 1.347849223658741
$dnArr[15,5] = 1.370630270025501
$dnArr[15,6] = 0.1905858059490519
$dnArr[15,7] = 0.78879439745441
$dnArr[15,8] = 0.3340804026300077
$dnArr[15,9] = 0.2531785942006834
$dnArr[15,10] = 2.522978054007925
$dnArr[16,0] = 0.09671376651664332
$dnArr[16,1] = 0.1284950338170014
$dnArr[16,2] = 2.131007511911506
$dnArr[16,3] = 1.48747951131287
$dnArr[16,4] = 1.348469572164731
$dnArr[16,5] = 1.372460363830172
$dnArr[16,6] = 0.1905208609297553
$dnArr[16,7] = 0.7742404952944071
$dnArr[16,8] = 0.3319567217420172
$dnArr[16,9] = 0.2514528043337734
$dnArr[16,10] = 2.528127343650095
$dnArr[17,0] = 0.09670154200265557
$dnArr[17,1] = 0.1285019797050972
$dnArr[17,2] = 2.130939934507779
$dnArr[17,3] = 1.487136713297573
$dnArr[17,4] = 1.348688583308856
$dnArr[17,5] = 1.373088260771503
$dnArr[17,6] = 0.1904989289952521
$dnArr[17,7] = 0.7693161499977919
$dnArr[17,8] = 0.3312400330619312
$dnArr[17,9] = 0.250870128911032
$dnArr[17,10] = 2.529882692480861
$dnArr[18,0] = 0.09675738729818306
$dnArr[18,1] = 0.1284712313628196
$dnArr[18,2] = 2.131307042435083
$dnArr[18,3] = 1.48873099644949
$dnArr[18,4] = 1.347738671973147
$dnArr[18,5] = 1.370295483054335
$dnArr[18,6] = 0.1905978530261194
$dnArr[18,7] = 0.7914895882175301
$dnArr[18,8] = 0.3344745651714049
$dnArr[18,9] = 0.2534987806956437
$dnArr[18,10] = 2.522030682918967
$dnArr[19,0] = 0.09695955253819477
$dnArr[19,1] = 0.1283768186827925
$dnArr[19,2] = 2.13364237358627
$dnArr[19,3] = 1.494991612895305
$dnArr[19,4] = 1.345170806961164
$dnArr[19,5] = 1.36148632136738
$dnArr[19,6] = 0.1909341999172334
$dnArr[19,7] = 0.8662205020283693
$dnArr[19,8] = 0.345505333250145
$dnArr[19,9] = 0.2624447151575353
$dnArr[19,10] = 2.49647604360376
$dnArr[20,0] = 0.09710244369725629
$dnArr[20,1] = 0.1283217356690785
$dnArr[20,2] = 2.135989131289477
$dnArr[20,3] = 1.499757199881913
$dnArr[20,4] = 1.343956932372052
$dnArr[20,5] = 1.356156601782494
$dnArr[20,6] = 0.191156970430864
$dnArr[20,7] = 0.9152311357546807
$dnArr[20,8] = 0.3528355950758026
$dnArr[20,9] = 0.2683757089459391
$dnArr[20,10] = 2.48039109628931
$dnArr[21,0] = 0.09702519930860376
$dnArr[21,1] = 0.1283505293248863
$dnArr[21,2] = 2.134661597948991
$dnArr[21,3] = 1.4971520458908
$dnArr[21,4] = 1.344562340772228
$dnArr[21,5] = 1.358962074018301
$dnArr[21,6] = 0.19103780477624
$dnArr[21,7] = 0.8890578846538233
$dnArr[21,8] = 0.3489122621603684
$dnArr[21,9] = 0.265202555074886
$dnArr[21,10] = 2.488919520868013
$dnArr[22,0] = 0.09675426569783596
$dnArr[22,1] = 0.1284728864156619
$dnArr[22,2] = 2.13128273525173
$dnArr[22,3] = 1.488640045283503
$dnArr[22,4] = 1.347788488914006
$dnArr[22,5] = 1.370446687984021
$dnArr[22,6] = 0.1905924055966839
$dnArr[22,7] = 0.7902710523959513
$dnArr[22,8] = 0.3342963246775383
$dnArr[22,9] = 0.2533539970126952
$dnArr[22,10] = 2.522458766613273
$dnArr[23,0] = 0.09650867053725776
$dnArr[23,1] = 0.1286333633412362
$dnArr[23,2] = 2.131165699576684
$dnArr[23,3] = 1.482350998233855
$dnArr[23,4] = 1.353265269628466
$dnArr[23,5] = 1.384672901178238
$dnArr[23,6] = 0.1901247233613717
$dnArr[23,7] = 0.684458664073702
$dnArr[23,8] = 0.3190541008548706
$dnArr[23,9] = 0.2409397035203824
$dnArr[23,10] = 2.561282666043546
$ws.Range("D2:N25").Value = $dnArr

Write-Host "Done updating pl_mw values"
